$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.932.51'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '1.633.17'
$ws.Range("E3").Value = '  -0.83%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("E5").Value = '  -0.67%  '

$ws.Range("E6").Value = '  -0.87%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.28'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("E9").Value = '  -2.73%  '

$ws.Range("E10").Value = '  -0.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0882'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.09%  '

$ws.Range("D12").Value = '1.864.78'
$ws.Range("E12").Value = '  -0.83%  '

$ws.Range("D13").Value = '1.629.50'
$ws.Range("E13").Value = '  -1.00%  '

$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.41%  '

$ws.Range("D17").Value = '27.931.74'
$ws.Range("E17").Value = '  -0.15%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '

$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.07%  '

$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.37'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("E24").Value = '  -3.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.83%  '

$ws.Range("E26").Value = '  +0.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("E28").Value = '  -0.54%  '

$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("E30").Value = '  -1.08%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0483'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("E32").Value = '  +0.99%  '

$ws.Range("B33").Value = 'Maker'
$ws.Range("C33").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D33").Value = '1.404.00'
$ws.Range("E33").Value = '  -2.82%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.11%  '

$ws.Range("E37").Value = '  +1.51%  '

$ws.Range("E38").Value = '  +0.77%  '

$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("E40").Value = '  -1.86%  '

$ws.Range("E41").Value = '  -0.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.42%  '

$ws.Range("E44").Value = '  +2.80%  '

$ws.Range("E45").Value = '  +1.25%  '

$ws.Range("E46").Value = '  -1.41%  '

$ws.Range("D47").Value = '1.774.29'
$ws.Range("E47").Value = '  -0.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.32%  '

$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  +0.76%  '

$ws.Range("E50").Value = '  -1.06%  '

$ws.Range("E51").Value = '  -0.25%  '

Write-Output "cryptos list updated"
